$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 1,128
$arr[0,0] = 44105
$arr[0,1] = 835339
$arr[0,2] = 2744
$arr[0,3] = 115691
$arr[0,4] = 67664
$arr[0,5] = 270466
$arr[0,6] = 29137
$arr[0,7] = 7267
$arr[0,8] = 6042
$arr[0,9] = 8500
$arr[0,10] = 9623
$arr[0,11] = 20313
$arr[0,12] = 4003
$arr[0,13] = 24037
$arr[0,14] = 33496
$arr[0,15] = 8161
$arr[0,16] = 11672
$arr[0,17] = 15384
$arr[0,18] = 15358
$arr[0,19] = 18531
$arr[0,20] = 16024
$arr[0,21] = 3803
$arr[0,22] = 3776
$arr[0,23] = 11137
$arr[0,24] = 31510
$arr[0,25] = 14115
$arr[0,26] = 12482
$arr[0,27] = 63267
$arr[0,28] = 2458
$arr[0,29] = 1478
$arr[0,30] = 770
$arr[0,31] = 478
$arr[0,32] = 929
$arr[0,33] = 533
$arr[0,34] = 785
$arr[0,35] = 2076
$arr[0,36] = 5966
$arr[0,37] = 38231
$arr[0,38] = 10280
$arr[0,39] = 2578
$arr[0,40] = 48737
$arr[0,41] = 1169
$arr[0,42] = 23243
$arr[0,43] = 1535
$arr[0,44] = 10599
$arr[0,45] = 1686
$arr[0,46] = 1616
$arr[0,47] = 8691
$arr[0,48] = 2063
$arr[0,49] = 966
$arr[0,50] = 2503
$arr[0,51] = 2700
$arr[0,52] = 66841
$arr[0,53] = 14255
$arr[0,54] = 6968
$arr[0,55] = 10037
$arr[0,56] = 7487
$arr[0,57] = 257
$arr[0,58] = 1467
$arr[0,59] = 2743
$arr[0,60] = 745
$arr[0,61] = 2177
$arr[0,62] = 10028
$arr[0,63] = 9604
$arr[0,64] = 10791
$arr[0,65] = 14395
$arr[0,66] = 1975
$arr[0,67] = 912
$arr[0,68] = 14337
$arr[0,69] = 11699
$arr[0,70] = 13735
$arr[0,71] = 3327
$arr[0,72] = 2329
$arr[0,73] = 6283
$arr[0,74] = 5077
$arr[0,75] = 2676
$arr[0,76] = 6275
$arr[0,77] = 3927
$arr[0,78] = 2314
$arr[0,79] = 1291
$arr[0,80] = 3085
$arr[0,81] = 2255
$arr[0,82] = 2144
$arr[0,83] = 2021
$arr[0,84] = 6619
$arr[0,85] = 2260
$arr[0,86] = 1533
$arr[0,87] = 1868
$arr[0,88] = 2154
$arr[0,89] = 2324
$arr[0,90] = 2698
$arr[0,91] = 1869
$arr[0,92] = 1231
$arr[0,93] = 1233
$arr[0,94] = 1142
$arr[0,95] = 3487
$arr[0,96] = 1523
$arr[0,97] = 985
$arr[0,98] = 1183
$arr[0,99] = 1767
$arr[0,100] = 1644
$arr[0,101] = 822
$arr[0,102] = 922
$arr[0,103] = 1397
$arr[0,104] = 1840
$arr[0,105] = 1724
$arr[0,106] = 1698
$arr[0,107] = 1305
$arr[0,108] = 335
$arr[0,109] = 373
$arr[0,110] = 841
$arr[0,111] = 790
$arr[0,112] = 509
$arr[0,113] = 544
$arr[0,114] = 392
$arr[0,115] = 698
$arr[0,116] = 763
$arr[0,117] = 532
$arr[0,118] = 508
$arr[0,119] = 374
$arr[0,120] = 527
$arr[0,121] = 140659
$arr[0,122] = 352375
$arr[0,123] = 21228
$arr[0,124] = 153833
$arr[0,125] = 95424
$arr[0,126] = 47583
$arr[0,127] = 13232
$ws.Range("A211:DX211").Value = $arr
[void]$ws.Range("CH189").Select()
